# Generate Report for Handback
# Update timestamps on the "Overview", "zh-cn" and "de-de" sheets to reflect
# the latest handback/xliff-generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to f26e92ad-55d0-40a4-bf24-14eee9e32772.md
$wsOverview.Range("G3").Value = "2016-08-12 06:56:20"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to f26e92ad-55d0-40a4-bf24-14eee9e32772...zh-cn.xlf
$wsZhCn.Range("H3").Value = "2016-08-12 06:56:13"
$wsZhCn.Range("K3").Value = "2016-08-12 06:56:44"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to f26e92ad-55d0-40a4-bf24-14eee9e32772...de-de.xlf
$wsDeDe.Range("H3").Value = "2016-08-12 06:56:20"
$wsDeDe.Range("K3").Value = "2016-08-12 06:56:53"
